$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "statut" column (A) encodes a status using one of four emoji glyphs:
#   📕 (red book)    -> rouge  (red)
#   📘 (blue book)   -> bleu   (blue)
#   📙 (orange book) -> orange (orange)
#   📗 (green book)  -> vert   (green)
# Excel renders these emoji inconsistently across platforms/fonts, so they
# are being replaced with plain, portable glyphs/text:
#   📕 -> -3
#   📘 -> ⚠️
#   📙 -> +3
#   📗 -> ✅

$lastRow = $ws.UsedRange.Rows.Count()

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $v = $cell.Value()

    if ($v -eq "📕") {
        # "-3" looks numeric to Excel's auto-detection, so force the cell
        # to Text first, write the value, then drop back to the default
        # style so no stray formatting is left behind on the cell.
        $cell.NumberFormat = "@"
        $cell.Value = "-3"
        $cell.Style = "Normal"
    } elseif ($v -eq "📘") {
        $cell.Value = "⚠️"
    } elseif ($v -eq "📙") {
        # "+3" is likewise auto-detected as a number by Excel.
        $cell.NumberFormat = "@"
        $cell.Value = "+3"
        $cell.Style = "Normal"
    } elseif ($v -eq "📗") {
        $cell.Value = "✅"
    }
}
